# Update time-series factor values for the social-care analysis "variant 2".
# Sheet "UK_cohabitation_adjustment" (B12:B61) and
# Sheet "UK_fertility_adjustment" (B12:B61) both get new, non-zero values.
# The cohabitation-adjustment sheet also becomes the active sheet/tab,
# matching the author's saved view state.

$wb = $excel.ActiveWorkbook

$cohabValues = @(
    -1.70385500364881,
    -1.1571997465369599,
    -1.80929478189082,
    -1.45740161376961,
    -1.29073044690233,
    -1.3056634031600001,
    -1.23057408899539,
    -1.22857101286559,
    -1.1826767128946001,
    -1.1579759478193301,
    -1.12963486475958,
    -1.0757812715078701,
    -1.12092793449733,
    -1.1018731590873201,
    -1.09760958793781,
    -1.0830768496758301,
    -1.0555063329510599,
    -1.0479160285220399,
    -1.0514706039277,
    -1.06908691107647,
    -1.0415116257308501,
    -1.0513753123935701,
    -1.04265085219326,
    -1.0414216342656399,
    -1.0324540674931599,
    -1.02569115746613,
    -1.02695584627324,
    -1.02822328009438,
    -1.0256289303811199,
    -1.01431235476579,
    -1.0189988164694199,
    -0.93569233126322005,
    -1.0965571763798101,
    -1.0328666908366799,
    -1.0312928178387699,
    -1.0409650808786699,
    -1.0051876412406699,
    -1.0203281201491501,
    -1.0024476094547801,
    -0.99635378684264198,
    -1.01498858055898,
    -1.0402525118708299,
    -0.99703465162276605,
    -1.0075997527413501,
    -0.96736979634708398,
    -0.95097258952846597,
    -0.94182162387448598,
    -0.99147068097634095,
    -0.93307533098083895,
    -0.94063451050118296
)

$fertilityValues = @(
    -0.46934212296623401,
    -0.454739741051573,
    -0.38374565123458099,
    -0.36783896829080098,
    -0.33165493609233199,
    -0.291562375123395,
    -0.25198600505337698,
    -0.28483236661986799,
    -0.19968626665699099,
    -0.18811664205453299,
    -0.19600519076745501,
    -0.170825426649971,
    -0.180453640042601,
    -0.16804309294200601,
    -0.17399727223845801,
    -0.17239504964128299,
    -0.190003699937489,
    -0.20430272643120501,
    -0.175676089303543,
    -0.20998728649817899,
    -0.203331600424143,
    -0.21884376646789699,
    -0.23494724441978199,
    -0.203708135524758,
    -0.21974647037595901,
    -0.22176718155054601,
    -0.216239660869704,
    -0.234189632538475,
    -0.20718286006263001,
    -0.24936479477594001,
    -0.22625695535786999,
    -0.22193325066712699,
    -0.255926102036931,
    -0.215144807846602,
    -0.205975543344144,
    -0.22623264483573799,
    -0.21071425092177301,
    -0.223422932715302,
    -0.241519153733185,
    -0.22311606421336,
    -0.247279492286858,
    -0.225909829302193,
    -0.21378817134198899,
    -0.22796885539012501,
    -0.23952823907995999,
    -0.23200492090432101,
    -0.25071493326128902,
    -0.22327190713982001,
    -0.22850366385618201,
    -0.21569328218244499
)

$wsCohab = $wb.Worksheets.Item("UK_cohabitation_adjustment")
for ($i = 0; $i -lt $cohabValues.Length; $i++) {
    $wsCohab.Cells.Item(12 + $i, 2).Value = $cohabValues[$i]
}

$wsFertility = $wb.Worksheets.Item("UK_fertility_adjustment")
for ($i = 0; $i -lt $fertilityValues.Length; $i++) {
    $wsFertility.Cells.Item(12 + $i, 2).Value = $fertilityValues[$i]
}

# Matches the workbook diff: UK_cohabitation_adjustment becomes the
# active/selected tab (activeTab 3 -> 2, tabSelected moves sheets).
$wsCohab.Activate()
